# Apply CryCompanywiseStockReport stock-count corrections.
# For each affected item row: Qty (F) is corrected and Value (G = UnitCost*Qty)
# is recomputed; a handful of rows had their Item Code/Rate/Qty/Value (B/E/F/G)
# swapped with the following row (two records that were out of order); and every
# "Sub Total:"/"Grand Total:" row (column B) is updated to the new sum of its section.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 70
$ws.Range("G6").Value = 2091.6

$ws.Range("B10").Value = 27352.87

$ws.Range("F21").Value = 148
$ws.Range("G21").Value = 3799.16

$ws.Range("B32").Value = 12602.26

$ws.Range("F50").Value = 1
$ws.Range("G50").Value = 629.0599999999999

$ws.Range("B52").Value = 5241.25

$ws.Range("F68").Value = 43
$ws.Range("G68").Value = 4950.16

$ws.Range("F70").Value = 8
$ws.Range("G70").Value = 1079.6

$ws.Range("F71").Value = 318
$ws.Range("G71").Value = 20256.6

$ws.Range("F77").Value = 248
$ws.Range("G77").Value = 11591.52

$ws.Range("F83").Value = 111
$ws.Range("G83").Value = 16724.37

$ws.Range("F84").Value = 26
$ws.Range("G84").Value = 2663.96

$ws.Range("F85").Value = 138
$ws.Range("G85").Value = 18598.26

$ws.Range("F86").Value = 58
$ws.Range("G86").Value = 7277.26

$ws.Range("F89").Value = 1
$ws.Range("G89").Value = 47.3

$ws.Range("B90").Value = 174179.97

$ws.Range("F102").Value = 5
$ws.Range("G102").Value = 247.4

$ws.Range("B104").Value = 215.38

$ws.Range("B112").Value = 57756
$ws.Range("E112").Value = 79.37
$ws.Range("F112").Value = -100
$ws.Range("G112").Value = -6644

$ws.Range("B113").Value = 64350
$ws.Range("E113").Value = 70.63
$ws.Range("F113").Value = 2
$ws.Range("G113").Value = 132.88

$ws.Range("F119").Value = 2
$ws.Range("G119").Value = 217.62

$ws.Range("B120").Value = 217.62

$ws.Range("B127").Value = 64329
$ws.Range("E127").Value = 128.32
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 120.69

$ws.Range("B128").Value = 57552
$ws.Range("E128").Value = 136.86
$ws.Range("F128").Value = -5
$ws.Range("G128").Value = -603.45

$ws.Range("F135").Value = 23
$ws.Range("G135").Value = 713.6900000000001

$ws.Range("B138").Value = 2544.83

$ws.Range("F141").Value = 47
$ws.Range("G141").Value = 2515.91

$ws.Range("B142").Value = 3027.59

$ws.Range("F144").Value = 1019
$ws.Range("G144").Value = 8610.549999999999

$ws.Range("F145").Value = 432
$ws.Range("G145").Value = 3451.68

$ws.Range("F146").Value = 20
$ws.Range("G146").Value = 1683.8

$ws.Range("B147").Value = 13746.03

$ws.Range("F151").Value = 91
$ws.Range("G151").Value = 7906.08

$ws.Range("B156").Value = 31056.75

$ws.Range("F169").Value = 2
$ws.Range("G169").Value = 287

$ws.Range("B175").Value = 26995.85

$ws.Range("F183").Value = 2
$ws.Range("G183").Value = 278.7

$ws.Range("F191").Value = 0
$ws.Range("G191").Value = 0

$ws.Range("B192").Value = 48706
$ws.Range("E192").Value = 39.8
$ws.Range("F192").Value = -144
$ws.Range("G192").Value = -4795.2

$ws.Range("B193").Value = 64973
$ws.Range("E193").Value = 35.4
$ws.Range("F193").Value = 2
$ws.Range("G193").Value = 66.59999999999999

$ws.Range("F203").Value = 55
$ws.Range("G203").Value = 1108.8

$ws.Range("F214").Value = 46
$ws.Range("G214").Value = 4034.2

$ws.Range("B216").Value = 37818.71

$ws.Range("F218").Value = 6
$ws.Range("G218").Value = 1297.32

$ws.Range("F225").Value = 75
$ws.Range("G225").Value = 8567.25

$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 65
$ws.Range("G227").Value = 9378.200000000001

$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32

$ws.Range("F229").Value = 58
$ws.Range("G229").Value = 8321.84

$ws.Range("B232").Value = 63510
$ws.Range("E232").Value = 50.66
$ws.Range("F232").Value = 115
$ws.Range("G232").Value = 5478.6

$ws.Range("B233").Value = 55356
$ws.Range("E233").Value = 54.04
$ws.Range("F233").Value = -158
$ws.Range("G233").Value = -7527.12

$ws.Range("F234").Value = 39
$ws.Range("G234").Value = 2001.48

$ws.Range("B243").Value = 63560
$ws.Range("E243").Value = 134.87
$ws.Range("F243").Value = 1
$ws.Range("G243").Value = 126.86

$ws.Range("B244").Value = 60325
$ws.Range("E244").Value = 151.57
$ws.Range("F244").Value = -102
$ws.Range("G244").Value = -12939.72

$ws.Range("F247").Value = 134
$ws.Range("G247").Value = 13923.94

$ws.Range("F255").Value = 541
$ws.Range("G255").Value = 92689.53

$ws.Range("F256").Value = 269
$ws.Range("G256").Value = 40664.73

$ws.Range("B260").Value = 179626.48

$ws.Range("F292").Value = 43
$ws.Range("G292").Value = 3580.61

$ws.Range("F302").Value = 43
$ws.Range("G302").Value = 9068.27

$ws.Range("F303").Value = 28
$ws.Range("G303").Value = 5904.92

$ws.Range("B304").Value = 171308.54

$ws.Range("F338").Value = 77
$ws.Range("G338").Value = 1824.9

$ws.Range("F345").Value = 47
$ws.Range("G345").Value = 2886.27

$ws.Range("B346").Value = 25502.1

$ws.Range("F354").Value = 14
$ws.Range("G354").Value = 960.26

$ws.Range("B358").Value = 34983.25

$ws.Range("B442").Value = 53319
$ws.Range("E442").Value = 310.64
$ws.Range("F442").Value = -6
$ws.Range("G442").Value = -1643.52

$ws.Range("B443").Value = 64810
$ws.Range("E443").Value = 291.22
$ws.Range("F443").Value = 4
$ws.Range("G443").Value = 1095.68

$ws.Range("F465").Value = 144
$ws.Range("G465").Value = 4727.52

$ws.Range("B475").Value = 45288.22

$ws.Range("F477").Value = 6
$ws.Range("G477").Value = 272.04

$ws.Range("B478").Value = 272.04

$ws.Range("F485").Value = 12
$ws.Range("G485").Value = 2105.64

$ws.Range("B488").Value = 29753.7

$ws.Range("F491").Value = 17
$ws.Range("G491").Value = 3026.34

$ws.Range("B493").Value = 11250.99

$ws.Range("F509").Value = 210
$ws.Range("G509").Value = 16879.8

$ws.Range("B510").Value = 22596.5

$ws.Range("F525").Value = 5
$ws.Range("G525").Value = 455.4

$ws.Range("B526").Value = 455.4

$ws.Range("B572").Value = 65362
$ws.Range("F572").Value = 20
$ws.Range("G572").Value = 817.4

$ws.Range("B573").Value = 65079
$ws.Range("F573").Value = 6
$ws.Range("G573").Value = 245.22

$ws.Range("F577").Value = 47
$ws.Range("G577").Value = 2020.53

$ws.Range("F578").Value = 68
$ws.Range("G578").Value = 3392.52

$ws.Range("F580").Value = 51
$ws.Range("G580").Value = 2906.49

$ws.Range("F582").Value = 25
$ws.Range("G582").Value = 1424.75

$ws.Range("B583").Value = 14295.75

$ws.Range("F599").Value = 1501
$ws.Range("G599").Value = 244828.11

$ws.Range("F601").Value = 386
$ws.Range("G601").Value = 109187.82

$ws.Range("F602").Value = 325
$ws.Range("G602").Value = 47011.25

$ws.Range("B606").Value = 401875.23

$ws.Range("F613").Value = 134
$ws.Range("G613").Value = 21327.44

$ws.Range("B618").Value = 43089.41

$ws.Range("B619").Value = 1674905.42

$ws.Range("B620").Value = 1674905.42
